# Scheduled-runner data refresh: update currentAveragePrice* / Leve* market
# figures (columns H-N) across the per-job Leve tables on each sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 680.3333
$ws.Range("I9").Value = 633
$ws.Range("K9").Value = 633
$ws.Range("M9").Value = -464
$ws.Range("H28").Value = 997.4666999999999
$ws.Range("I28").Value = 888.46155
$ws.Range("J28").Value = 1706
$ws.Range("K28").Value = 888.46155
$ws.Range("L28").Value = 1706
$ws.Range("M28").Value = -403.46155
$ws.Range("N28").Value = -2676
$ws.Range("H96").Value = 7937174
$ws.Range("I96").Value = 8928946
$ws.Range("K96").Value = 26786838
$ws.Range("M96").Value = -26785465
$ws.Range("H104").Value = 183
$ws.Range("I104").Value = 183
$ws.Range("K104").Value = 549
$ws.Range("M104").Value = 1198
$ws.Range("H111").Value = 1250
$ws.Range("H138").Value = 2328.5
$ws.Range("J138").Value = 4596.05
$ws.Range("L138").Value = 13788.15
$ws.Range("N138").Value = -24068.15
$ws.Range("H141").Value = 1840
$ws.Range("I141").Value = 1775.1818
$ws.Range("K141").Value = 5325.5454
$ws.Range("M141").Value = -145.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32688.725
$ws.Range("I32").Value = 24598.553
$ws.Range("J32").Value = 127748.25
$ws.Range("K32").Value = 24598.553
$ws.Range("L32").Value = 127748.25
$ws.Range("M32").Value = -24311.553
$ws.Range("N32").Value = -128322.25
$ws.Range("H61").Value = 4278.125
$ws.Range("I61").Value = 1485.2273
$ws.Range("J61").Value = 7691.6665
$ws.Range("K61").Value = 1485.2273
$ws.Range("L61").Value = 7691.6665
$ws.Range("M61").Value = -1273.2273
$ws.Range("N61").Value = -8115.6665
$ws.Range("H63").Value = 2916.5
$ws.Range("I63").Value = 2400
$ws.Range("K63").Value = 2400
$ws.Range("M63").Value = -1714
$ws.Range("H66").Value = 2916.5
$ws.Range("I66").Value = 2400
$ws.Range("K66").Value = 12000
$ws.Range("M66").Value = -8568
$ws.Range("H74").Value = 278264.9
$ws.Range("J74").Value = 34333.332
$ws.Range("L74").Value = 34333.332
$ws.Range("N74").Value = -36081.332
$ws.Range("H77").Value = 278264.9
$ws.Range("J77").Value = 34333.332
$ws.Range("L77").Value = 171666.66
$ws.Range("N77").Value = -180402.66
$ws.Range("H102").Value = 7545.2
$ws.Range("I102").Value = 6352.5
$ws.Range("K102").Value = 6352.5
$ws.Range("M102").Value = -4730.5
$ws.Range("H110").Value = 1532.1111
$ws.Range("I110").Value = 1020.61536
$ws.Range("K110").Value = 1020.61536
$ws.Range("M110").Value = 1024.38464
$ws.Range("H132").Value = 1363.262
$ws.Range("I132").Value = 1052.3704
$ws.Range("K132").Value = 3157.1112
$ws.Range("M132").Value = -627.1112000000003
$ws.Range("H136").Value = 4278.125
$ws.Range("I136").Value = 1485.2273
$ws.Range("J136").Value = 7691.6665
$ws.Range("K136").Value = 4455.6819
$ws.Range("L136").Value = 23074.9995
$ws.Range("M136").Value = -1905.6819
$ws.Range("N136").Value = -28174.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 87590.836
$ws.Range("I107").Value = 126108.25
$ws.Range("J107").Value = 10556
$ws.Range("K107").Value = 126108.25
$ws.Range("L107").Value = 10556
$ws.Range("M107").Value = -124188.25
$ws.Range("N107").Value = -14396
$ws.Range("H134").Value = 1728.1404
$ws.Range("I134").Value = 1335.2291
$ws.Range("J134").Value = 3823.6667
$ws.Range("K134").Value = 4005.6873
$ws.Range("L134").Value = 11471.0001
$ws.Range("M134").Value = -1470.6873
$ws.Range("N134").Value = -16541.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 74125
$ws.Range("I86").Value = 113731.164
$ws.Range("J86").Value = 26597.6
$ws.Range("K86").Value = 113731.164
$ws.Range("L86").Value = 26597.6
$ws.Range("M86").Value = -112608.164
$ws.Range("N86").Value = -28843.6
$ws.Range("H89").Value = 74125
$ws.Range("I89").Value = 113731.164
$ws.Range("J89").Value = 26597.6
$ws.Range("K89").Value = 568655.8200000001
$ws.Range("L89").Value = 132988
$ws.Range("M89").Value = -563039.8200000001
$ws.Range("N89").Value = -144220
$ws.Range("H99").Value = 6952.421
$ws.Range("I99").Value = 5674.467
$ws.Range("K99").Value = 5674.467
$ws.Range("M99").Value = -4176.467
$ws.Range("H126").Value = 6952.421
$ws.Range("I126").Value = 5674.467
$ws.Range("K126").Value = 17023.401
$ws.Range("M126").Value = -14553.401
$ws.Range("H132").Value = 43058.957
$ws.Range("I132").Value = 56493.11
$ws.Range("J132").Value = 2756.5
$ws.Range("K132").Value = 169479.33
$ws.Range("L132").Value = 8269.5
$ws.Range("M132").Value = -166949.33
$ws.Range("N132").Value = -13329.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 222.54546
$ws.Range("I8").Value = 222.54546
$ws.Range("K8").Value = 667.6363799999999
$ws.Range("M8").Value = -528.6363799999999
$ws.Range("H62").Value = 4868.421
$ws.Range("H65").Value = 4868.421
$ws.Range("H113").Value = 1254.0667
$ws.Range("I113").Value = 399.2
$ws.Range("J113").Value = 1681.5
$ws.Range("K113").Value = 1197.6
$ws.Range("L113").Value = 5044.5
$ws.Range("M113").Value = 972.4000000000001
$ws.Range("N113").Value = -9384.5
$ws.Range("H115").Value = 656.25
$ws.Range("I115").Value = 208.33333
$ws.Range("K115").Value = 624.99999
$ws.Range("M115").Value = 550.00001
$ws.Range("H131").Value = 1852.44
$ws.Range("J131").Value = 2119.2666
$ws.Range("L131").Value = 6357.7998
$ws.Range("N131").Value = -16437.7998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 40001.668
$ws.Range("J49").Value = 40001.668
$ws.Range("L49").Value = 40001.668
$ws.Range("N49").Value = -40369.668
$ws.Range("H80").Value = 7806.25
$ws.Range("I80").Value = 5140.6
$ws.Range("K80").Value = 5140.6
$ws.Range("M80").Value = -4142.6
$ws.Range("H83").Value = 7806.25
$ws.Range("I83").Value = 5140.6
$ws.Range("K83").Value = 25703
$ws.Range("M83").Value = -20711
$ws.Range("H102").Value = 51492.5
$ws.Range("I102").Value = 55901
$ws.Range("J102").Value = 2999
$ws.Range("K102").Value = 55901
$ws.Range("L102").Value = 2999
$ws.Range("M102").Value = -54279
$ws.Range("N102").Value = -6243
$ws.Range("H132").Value = 3923.6428
$ws.Range("I132").Value = 3869.375
$ws.Range("J132").Value = 4249.25
$ws.Range("K132").Value = 11608.125
$ws.Range("L132").Value = 12747.75
$ws.Range("M132").Value = -9078.125
$ws.Range("N132").Value = -17807.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3272.5557
$ws.Range("I16").Value = 3292.8
$ws.Range("K16").Value = 3292.8
$ws.Range("M16").Value = -3122.8
$ws.Range("H22").Value = 1146.1111
$ws.Range("I22").Value = 528.3333
$ws.Range("J22").Value = 1455
$ws.Range("K22").Value = 528.3333
$ws.Range("L22").Value = 1455
$ws.Range("M22").Value = -233.3333
$ws.Range("N22").Value = -2045
$ws.Range("H27").Value = 1146.1111
$ws.Range("I27").Value = 528.3333
$ws.Range("J27").Value = 1455
$ws.Range("K27").Value = 528.3333
$ws.Range("L27").Value = 1455
$ws.Range("M27").Value = -421.3333
$ws.Range("N27").Value = -1669
$ws.Range("H61").Value = 960.6667
$ws.Range("I61").Value = 1041
$ws.Range("J61").Value = 800
$ws.Range("K61").Value = 1041
$ws.Range("L61").Value = 800
$ws.Range("M61").Value = -839
$ws.Range("N61").Value = -1204
$ws.Range("H82").Value = 1912.5714
$ws.Range("J82").Value = 1414.6666
$ws.Range("L82").Value = 1414.6666
$ws.Range("N82").Value = -2136.6666
$ws.Range("H85").Value = 1912.5714
$ws.Range("J85").Value = 1414.6666
$ws.Range("L85").Value = 1414.6666
$ws.Range("N85").Value = -3910.6666
$ws.Range("H92").Value = 68330
$ws.Range("J92").Value = 68330
$ws.Range("L92").Value = 68330
$ws.Range("N92").Value = -73322
$ws.Range("H113").Value = 960.6667
$ws.Range("I113").Value = 1041
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 1041
$ws.Range("L113").Value = 800
$ws.Range("M113").Value = 1129
$ws.Range("N113").Value = -5140
$ws.Range("H132").Value = 1971.0303
$ws.Range("I132").Value = 1567.84
$ws.Range("K132").Value = 4703.52
$ws.Range("M132").Value = -2173.52
$ws.Range("H136").Value = 2727.64
$ws.Range("J136").Value = 2800.5715
$ws.Range("L136").Value = 8401.7145
$ws.Range("N136").Value = -13501.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 50000
$ws.Range("J58").Value = 40000
$ws.Range("L58").Value = 40000
$ws.Range("N58").Value = -40616
$ws.Range("H122").Value = 62917.684
$ws.Range("I122").Value = 96188.71000000001
$ws.Range("K122").Value = 288566.13
$ws.Range("M122").Value = -286116.13
$ws.Range("H132").Value = 2901.8774
$ws.Range("I132").Value = 2551.0466
$ws.Range("J132").Value = 5416.1665
$ws.Range("K132").Value = 7653.139800000001
$ws.Range("L132").Value = 16248.4995
$ws.Range("M132").Value = -5123.139800000001
$ws.Range("N132").Value = -21308.4995
$ws.Range("H136").Value = 13616
$ws.Range("I136").Value = 14908.609
$ws.Range("K136").Value = 44725.827
$ws.Range("M136").Value = -42175.827
